# Mark the remaining "ACTIVITY DIAGRAM" checklist rows as done by putting
# an "x" in column C, and move the sheet's view/selection down to the end
# of the list (mirrors a user having just finished filling these rows in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 48, 49, 52, 54, 55, 56, 57
foreach ($r in $rows) {
    $ws.Range("C$r").Value = "x"
}

$ws.Application.ActiveWindow.ScrollRow = 36
$ws.Range("C57").Select()
